{"js": "const replacements = [\n  [\"2024-12-21 Saturday\", \"2024-12-22 Sunday\"],\n  [\"82-10=72\", \"14+11=25\"],\n  [\"12-6=6\", \"50-32=18\"],\n  [\"48-21=27\", \"96-77=19\"],\n  [\"11+38=49\", \"24+74=98\"],\n  [\"34+50=84\", \"72-44=28\"],\n  [\"88-82=6\", \"54-30=24\"],\n  [\"4+16=20\", \"74-25=49\"],\n  [\"32-2=30\", \"25+18=43\"],\n  [\"12+8=20\", \"91-11=80\"],\n  [\"58-54=4\", \"95-51=44\"],\n  [\"50+35=85\", \"8+38=46\"],\n  [\"66+16=82\", \"50+49=99\"],\n  [\"62+11=73\", \"45+22=67\"],\n  [\"23+11=34\", \"54+40=94\"],\n  [\"83-28=55\", \"8+24=32\"],\n  [\"18-9=9\", \"22+16=38\"],\n  [\"77-14=63\", \"47+43=90\"],\n  [\"10+61=71\", \"2+17=19\"],\n  [\"83-16=67\", \"87-33=54\"],\n  [\"61+18=79\", \"30+0=30\"],\n  [\"12+84=96\", \"58+2=60\"],\n  [\"35+45=80\", \"37+43=80\"],\n  [\"74-41=33\", \"33-19=14\"],\n  [\"15-11=4\", \"27+56=83\"],\n  [\"71-58=13\", \"67-49=18\"],\n  [\"7+75=82\", \"63-1=62\"],\n  [\"31-16=15\", \"68+14=82\"],\n  [\"33+10=43\", \"33+39=72\"],\n  [\"39-6=33\", \"66-21=45\"],\n  [\"93-77=16\", \"87-69=18\"],\n  [\"44+33=77\", \"47+44=91\"],\n  [\"45+44=89\", \"77-28=49\"],\n  [\"42+36=78\", \"17+66=83\"],\n  [\"40+31=71\", \"39+5=44\"],\n  [\"69-57=12\", \"12+35=47\"],\n  [\"15+76=91\", \"70-34=36\"],\n  [\"65-57=8\", \"86-6=80\"],\n  [\"12+0=12\", \"25-11=14\"],\n  [\"85-64=21\", \"34+18=52\"],\n  [\"34+33=67\", \"93-90=3\"],\n  [\"45-12=33\", \"43-14=29\"],\n  [\"22+34=56\", \"63-15=48\"],\n  [\"79-78=1\", \"24+40=64\"],\n  [\"11+67=78\", \"46-3=43\"],\n  [\"84-9=75\", \"82-61=21\"],\n  [\"45+39=84\", \"57-32=25\"],\n  [\"46+44=90\", \"80+1=81\"],\n  [\"16+49=65\", \"72-7=65\"],\n  [\"63+21=84\", \"57-1=56\"],\n  [\"57+3=60\", \"66+24=90\"],\n  [\"7+10=17\", \"2+16=18\"],\n  [\"75-52=23\", \"23-8=15\"],\n  [\"96-18=78\", \"51+19=70\"],\n  [\"91-7=84\", \"37+20=57\"],\n  [\"11+49=60\", \"81-48=33\"],\n  [\"23-16=7\", \"99-19=80\"],\n  [\"63-20=43\", \"38+18=56\"],\n  [\"65-28=37\", \"70+9=79\"],\n  [\"25+54=79\", \"32+6=38\"],\n  [\"72+3=75\", \"78-62=16\"],\n  [\"47-2=45\", \"48+13=61\"],\n  [\"36-9=27\", \"20+26=46\"],\n  [\"45-7=38\", \"68+19=87\"],\n  [\"98-2=96\", \"83+15=98\"],\n  [\"18+73=91\", \"41+44=85\"],\n  [\"41+37=78\", \"43-27=16\"],\n  [\"56-39=17\", \"28-22=6\"],\n  [\"4+77=81\", \"46+5=51\"],\n  [\"27+26=53\", \"40+53=93\"],\n  [\"56+36=92\", \"77-51=26\"],\n  [\"6+63=69\", \"64-53=11\"],\n  [\"13+31=44\", \"11+81=92\"],\n  [\"82-1=81\", \"49-24=25\"],\n  [\"40+32=72\", \"16+52=68\"],\n  [\"40+29=69\", \"20+75=95\"],\n  [\"88-6=82\", \"65-10=55\"],\n  [\"76-8=68\", \"50+15=65\"],\n  [\"83-67=16\", \"43+47=90\"],\n  [\"62-22=40\", \"35-3=32\"],\n  [\"18-15=3\", \"40-31=9\"],\n  [\"62-31=31\", \"74-18=56\"],\n  [\"90-34=56\", \"13+20=33\"],\n  [\"14+67=81\", \"41+5=46\"],\n  [\"43-34=9\", \"37+25=62\"],\n  [\"53-13=40\", \"27+9=36\"],\n  [\"16+4=20\", \"32+14=46\"],\n  [\"93-66=27\", \"90-58=32\"],\n  [\"10+47=57\", \"63+31=94\"],\n  [\"20+21=41\", \"77-64=13\"],\n  [\"39-24=15\", \"70-17=53\"],\n  [\"38+22=60\", \"41+52=93\"],\n  [\"93-79=14\", \"92+6=98\"],\n  [\"44+51=95\", \"94-43=51\"],\n  [\"48-22=26\", \"26+19=45\"],\n  [\"14+16=30\", \"83-58=25\"],\n  [\"59+10=69\", \"39+50=89\"],\n  [\"37-32=5\", \"16+6=22\"],\n  [\"21+1=22\", \"82-28=54\"],\n  [\"13+82=95\", \"62+31=93\"],\n  [\"45+48=93\", \"2+80=82\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$null = $d.Content.Find.Execute(\"2024-12-21 Saturday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-12-22 Sunday\", 2)\n$null = $d.Content.Find.Execute(\"82-10=72\", $false, $false, $false, $false, $false, $true, 1, $false, \"14+11=25\", 2)\n$null = $d.Content.Find.Execute(\"12-6=6\", $false, $false, $false, $false, $false, $true, 1, $false, \"50-32=18\", 2)\n$null = $d.Content.Find.Execute(\"48-21=27\", $false, $false, $false, $false, $false, $true, 1, $false, \"96-77=19\", 2)\n$null = $d.Content.Find.Execute(\"11+38=49\", $false, $false, $false, $false, $false, $true, 1, $false, \"24+74=98\", 2)\n$null = $d.Content.Find.Execute(\"34+50=84\", $false, $false, $false, $false, $false, $true, 1, $false, \"72-44=28\", 2)\n$null = $d.Content.Find.Execute(\"88-82=6\", $false, $false, $false, $false, $false, $true, 1, $false, \"54-30=24\", 2)\n$null = $d.Content.Find.Execute(\"4+16=20\", $false, $false, $false, $false, $false, $true, 1, $false, \"74-25=49\", 2)\n$null = $d.Content.Find.Execute(\"32-2=30\", $false, $false, $false, $false, $false, $true, 1, $false, \"25+18=43\", 2)\n$null = $d.Content.Find.Execute(\"12+8=20\", $false, $false, $false, $false, $false, $true, 1, $false, \"91-11=80\", 2)\n$null = $d.Content.Find.Execute(\"58-54=4\", $false, $false, $false, $false, $false, $true, 1, $false, \"95-51=44\", 2)\n$null = $d.Content.Find.Execute(\"50+35=85\", $false, $false, $false, $false, $false, $true, 1, $false, \"8+38=46\", 2)\n$null = $d.Content.Find.Execute(\"66+16=82\", $false, $false, $false, $false, $false, $true, 1, $false, \"50+49=99\", 2)\n$null = $d.Content.Find.Execute(\"62+11=73\", $false, $false, $false, $false, $false, $true, 1, $false, \"45+22=67\", 2)\n$null = $d.Content.Find.Execute(\"23+11=34\", $false, $false, $false, $false, $false, $true, 1, $false, \"54+40=94\", 2)\n$null = $d.Content.Find.Execute(\"83-28=55\", $false, $false, $false, $false, $false, $true, 1, $false, \"8+24=32\", 2)\n$null = $d.Content.Find.Execute(\"18-9=9\", $false, $false, $false, $false, $false, $true, 1, $false, \"22+16=38\", 2)\n$null = $d.Content.Find.Execute(\"77-14=63\", $false, $false, $false, $false, $false, $true, 1, $false, \"47+43=90\", 2)\n$null = $d.Content.Find.Execute(\"10+61=71\", $false, $false, $false, $false, $false, $true, 1, $false, \"2+17=19\", 2)\n$null = $d.Content.Find.Execute(\"83-16=67\", $false, $false, $false, $false, $false, $true, 1, $false, \"87-33=54\", 2)\n$null = $d.Content.Find.Execute(\"61+18=79\", $false, $false, $false, $false, $false, $true, 1, $false, \"30+0=30\", 2)\n$null = $d.Content.Find.Execute(\"12+84=96\", $false, $false, $false, $false, $false, $true, 1, $false, \"58+2=60\", 2)\n$null = $d.Content.Find.Execute(\"35+45=80\", $false, $false, $false, $false, $false, $true, 1, $false, \"37+43=80\", 2)\n$null = $d.Content.Find.Execute(\"74-41=33\", $false, $false, $false, $false, $false, $true, 1, $false, \"33-19=14\", 2)\n$null = $d.Content.Find.Execute(\"15-11=4\", $false, $false, $false, $false, $false, $true, 1, $false, \"27+56=83\", 2)\n$null = $d.Content.Find.Execute(\"71-58=13\", $false, $false, $false, $false, $false, $true, 1, $false, \"67-49=18\", 2)\n$null = $d.Content.Find.Execute(\"7+75=82\", $false, $false, $false, $false, $false, $true, 1, $false, \"63-1=62\", 2)\n$null = $d.Content.Find.Execute(\"31-16=15\", $false, $false, $false, $false, $false, $true, 1, $false, \"68+14=82\", 2)\n$null = $d.Content.Find.Execute(\"33+10=43\", $false, $false, $false, $false, $false, $true, 1, $false, \"33+39=72\", 2)\n$null = $d.Content.Find.Execute(\"39-6=33\", $false, $false, $false, $false, $false, $true, 1, $false, \"66-21=45\", 2)\n$null = $d.Content.Find.Execute(\"93-77=16\", $false, $false, $false, $false, $false, $true, 1, $false, \"87-69=18\", 2)\n$null = $d.Content.Find.Execute(\"44+33=77\", $false, $false, $false, $false, $false, $true, 1, $false, \"47+44=91\", 2)\n$null = $d.Content.Find.Execute(\"45+44=89\", $false, $false, $false, $false, $false, $true, 1, $false, \"77-28=49\", 2)\n$null = $d.Content.Find.Execute(\"42+36=78\", $false, $false, $false, $false, $false, $true, 1, $false, \"17+66=83\", 2)\n$null = $d.Content.Find.Execute(\"40+31=71\", $false, $false, $false, $false, $false, $true, 1, $false, \"39+5=44\", 2)\n$null = $d.Content.Find.Execute(\"69-57=12\", $false, $false, $false, $false, $false, $true, 1, $false, \"12+35=47\", 2)\n$null = $d.Content.Find.Execute(\"15+76=91\", $false, $false, $false, $false, $false, $true, 1, $false, \"70-34=36\", 2)\n$null = $d.Content.Find.Execute(\"65-57=8\", $false, $false, $false, $false, $false, $true, 1, $false, \"86-6=80\", 2)\n$null = $d.Content.Find.Execute(\"12+0=12\", $false, $false, $false, $false, $false, $true, 1, $false, \"25-11=14\", 2)\n$null = $d.Content.Find.Execute(\"85-64=21\", $false, $false, $false, $false, $false, $true, 1, $false, \"34+18=52\", 2)\n$null = $d.Content.Find.Execute(\"34+33=67\", $false, $false, $false, $false, $false, $true, 1, $false, \"93-90=3\", 2)\n$null = $d.Content.Find.Execute(\"45-12=33\", $false, $false, $false, $false, $false, $true, 1, $false, \"43-14=29\", 2)\n$null = $d.Content.Find.Execute(\"22+34=56\", $false, $false, $false, $false, $false, $true, 1, $false, \"63-15=48\", 2)\n$null = $d.Content.Find.Execute(\"79-78=1\", $false, $false, $false, $false, $false, $true, 1, $false, \"24+40=64\", 2)\n$null = $d.Content.Find.Execute(\"11+67=78\", $false, $false, $false, $false, $false, $true, 1, $false, \"46-3=43\", 2)\n$null = $d.Content.Find.Execute(\"84-9=75\", $false, $false, $false, $false, $false, $true, 1, $false, \"82-61=21\", 2)\n$null = $d.Content.Find.Execute(\"45+39=84\", $false, $false, $false, $false, $false, $true, 1, $false, \"57-32=25\", 2)\n$null = $d.Content.Find.Execute(\"46+44=90\", $false, $false, $false, $false, $false, $true, 1, $false, \"80+1=81\", 2)\n$null = $d.Content.Find.Execute(\"16+49=65\", $false, $false, $false, $false, $false, $true, 1, $false, \"72-7=65\", 2)\n$null = $d.Content.Find.Execute(\"63+21=84\", $false, $false, $false, $false, $false, $true, 1, $false, \"57-1=56\", 2)\n$null = $d.Content.Find.Execute(\"57+3=60\", $false, $false, $false, $false, $false, $true, 1, $false, \"66+24=90\", 2)\n$null = $d.Content.Find.Execute(\"7+10=17\", $false, $false, $false, $false, $false, $true, 1, $false, \"2+16=18\", 2)\n$null = $d.Content.Find.Execute(\"75-52=23\", $false, $false, $false, $false, $false, $true, 1, $false, \"23-8=15\", 2)\n$null = $d.Content.Find.Execute(\"96-18=78\", $false, $false, $false, $false, $false, $true, 1, $false, \"51+19=70\", 2)\n$null = $d.Content.Find.Execute(\"91-7=84\", $false, $false, $false, $false, $false, $true, 1, $false, \"37+20=57\", 2)\n$null = $d.Content.Find.Execute(\"11+49=60\", $false, $false, $false, $false, $false, $true, 1, $false, \"81-48=33\", 2)\n$null = $d.Content.Find.Execute(\"23-16=7\", $false, $false, $false, $false, $false, $true, 1, $false, \"99-19=80\", 2)\n$null = $d.Content.Find.Execute(\"63-20=43\", $false, $false, $false, $false, $false, $true, 1, $false, \"38+18=56\", 2)\n$null = $d.Content.Find.Execute(\"65-28=37\", $false, $false, $false, $false, $false, $true, 1, $false, \"70+9=79\", 2)\n$null = $d.Content.Find.Execute(\"25+54=79\", $false, $false, $false, $false, $false, $true, 1, $false, \"32+6=38\", 2)\n$null = $d.Content.Find.Execute(\"72+3=75\", $false, $false, $false, $false, $false, $true, 1, $false, \"78-62=16\", 2)\n$null = $d.Content.Find.Execute(\"47-2=45\", $false, $false, $false, $false, $false, $true, 1, $false, \"48+13=61\", 2)\n$null = $d.Content.Find.Execute(\"36-9=27\", $false, $false, $false, $false, $false, $true, 1, $false, \"20+26=46\", 2)\n$null = $d.Content.Find.Execute(\"45-7=38\", $false, $false, $false, $false, $false, $true, 1, $false, \"68+19=87\", 2)\n$null = $d.Content.Find.Execute(\"98-2=96\", $false, $false, $false, $false, $false, $true, 1, $false, \"83+15=98\", 2)\n$null = $d.Content.Find.Execute(\"18+73=91\", $false, $false, $false, $false, $false, $true, 1, $false, \"41+44=85\", 2)\n$null = $d.Content.Find.Execute(\"41+37=78\", $false, $false, $false, $false, $false, $true, 1, $false, \"43-27=16\", 2)\n$null = $d.Content.Find.Execute(\"56-39=17\", $false, $false, $false, $false, $false, $true, 1, $false, \"28-22=6\", 2)\n$null = $d.Content.Find.Execute(\"4+77=81\", $false, $false, $false, $false, $false, $true, 1, $false, \"46+5=51\", 2)\n$null = $d.Content.Find.Execute(\"27+26=53\", $false, $false, $false, $false, $false, $true, 1, $false, \"40+53=93\", 2)\n$null = $d.Content.Find.Execute(\"56+36=92\", $false, $false, $false, $false, $false, $true, 1, $false, \"77-51=26\", 2)\n$null = $d.Content.Find.Execute(\"6+63=69\", $false, $false, $false, $false, $false, $true, 1, $false, \"64-53=11\", 2)\n$null = $d.Content.Find.Execute(\"13+31=44\", $false, $false, $false, $false, $false, $true, 1, $false, \"11+81=92\", 2)\n$null = $d.Content.Find.Execute(\"82-1=81\", $false, $false, $false, $false, $false, $true, 1, $false, \"49-24=25\", 2)\n$null = $d.Content.Find.Execute(\"40+32=72\", $false, $false, $false, $false, $false, $true, 1, $false, \"16+52=68\", 2)\n$null = $d.Content.Find.Execute(\"40+29=69\", $false, $false, $false, $false, $false, $true, 1, $false, \"20+75=95\", 2)\n$null = $d.Content.Find.Execute(\"88-6=82\", $false, $false, $false, $false, $false, $true, 1, $false, \"65-10=55\", 2)\n$null = $d.Content.Find.Execute(\"76-8=68\", $false, $false, $false, $false, $false, $true, 1, $false, \"50+15=65\", 2)\n$null = $d.Content.Find.Execute(\"83-67=16\", $false, $false, $false, $false, $false, $true, 1, $false, \"43+47=90\", 2)\n$null = $d.Content.Find.Execute(\"62-22=40\", $false, $false, $false, $false, $false, $true, 1, $false, \"35-3=32\", 2)\n$null = $d.Content.Find.Execute(\"18-15=3\", $false, $false, $false, $false, $false, $true, 1, $false, \"40-31=9\", 2)\n$null = $d.Content.Find.Execute(\"62-31=31\", $false, $false, $false, $false, $false, $true, 1, $false, \"74-18=56\", 2)\n$null = $d.Content.Find.Execute(\"90-34=56\", $false, $false, $false, $false, $false, $true, 1, $false, \"13+20=33\", 2)\n$null = $d.Content.Find.Execute(\"14+67=81\", $false, $false, $false, $false, $false, $true, 1, $false, \"41+5=46\", 2)\n$null = $d.Content.Find.Execute(\"43-34=9\", $false, $false, $false, $false, $false, $true, 1, $false, \"37+25=62\", 2)\n$null = $d.Content.Find.Execute(\"53-13=40\", $false, $false, $false, $false, $false, $true, 1, $false, \"27+9=36\", 2)\n$null = $d.Content.Find.Execute(\"16+4=20\", $false, $false, $false, $false, $false, $true, 1, $false, \"32+14=46\", 2)\n$null = $d.Content.Find.Execute(\"93-66=27\", $false, $false, $false, $false, $false, $true, 1, $false, \"90-58=32\", 2)\n$null = $d.Content.Find.Execute(\"10+47=57\", $false, $false, $false, $false, $false, $true, 1, $false, \"63+31=94\", 2)\n$null = $d.Content.Find.Execute(\"20+21=41\", $false, $false, $false, $false, $false, $true, 1, $false, \"77-64=13\", 2)\n$null = $d.Content.Find.Execute(\"39-24=15\", $false, $false, $false, $false, $false, $true, 1, $false, \"70-17=53\", 2)\n$null = $d.Content.Find.Execute(\"38+22=60\", $false, $false, $false, $false, $false, $true, 1, $false, \"41+52=93\", 2)\n$null = $d.Content.Find.Execute(\"93-79=14\", $false, $false, $false, $false, $false, $true, 1, $false, \"92+6=98\", 2)\n$null = $d.Content.Find.Execute(\"44+51=95\", $false, $false, $false, $false, $false, $true, 1, $false, \"94-43=51\", 2)\n$null = $d.Content.Find.Execute(\"48-22=26\", $false, $false, $false, $false, $false, $true, 1, $false, \"26+19=45\", 2)\n$null = $d.Content.Find.Execute(\"14+16=30\", $false, $false, $false, $false, $false, $true, 1, $false, \"83-58=25\", 2)\n$null = $d.Content.Find.Execute(\"59+10=69\", $false, $false, $false, $false, $false, $true, 1, $false, \"39+50=89\", 2)\n$null = $d.Content.Find.Execute(\"37-32=5\", $false, $false, $false, $false, $false, $true, 1, $false, \"16+6=22\", 2)\n$null = $d.Content.Find.Execute(\"21+1=22\", $false, $false, $false, $false, $false, $true, 1, $false, \"82-28=54\", 2)\n$null = $d.Content.Find.Execute(\"13+82=95\", $false, $false, $false, $false, $false, $true, 1, $false, \"62+31=93\", 2)\n$null = $d.Content.Find.Execute(\"45+48=93\", $false, $false, $false, $false, $false, $true, 1, $false, \"2+80=82\", 2)\n"}
